$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value = "24/10/2025"
$ws.Range("B9").Value = "Bremen"
$ws.Range("C9").Value = 1
$ws.Range("D9").Value = 0
$ws.Range("E9").Value = "Union Berlin"
$ws.Range("F9").Value = "W"
$ws.Range("G9").Value = 0
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 1
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 1.02
$ws.Range("L9").Value = 0.57
$ws.Range("M9").Value = 16
$ws.Range("N9").Value = 12
$ws.Range("O9").Value = 3
$ws.Range("P9").Value = 2
